$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Report generated timestamp
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"

# ---------------------------------------------------------------------------
# 2. Summary totals (Total Billed Amount / Total Line Items)
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 22

# ---------------------------------------------------------------------------
# 3. Monday (07/14/2025) block — edits to existing rows (before the new row
#    is inserted at row 25, so these row numbers are still the "before" ones)
# ---------------------------------------------------------------------------
$ws.Range("H16").Value = 0

# Row 17 changes identity from "Point 07 / PLA-HDIG" to "Point 11 / GND-MD"
$ws.Range("A17").Value = "Point 11"
$ws.Range("B17").Value = "GND-MD"
$ws.Range("D17").Value = "GND,Wire Mldg Only"
$ws.Range("F17").Value = 2
$ws.Range("H17").Value = 0

$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0

# Row 20 changes identity label from "Point 09" to "Point 07"
$ws.Range("A20").Value = "Point 07"
$ws.Range("H20").Value = 0

$ws.Range("H21").Value = 0

# Row 22 changes identity label from "Point 11" to "Point 09"
$ws.Range("A22").Value = "Point 09"
$ws.Range("H22").Value = 0

$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 0

# ---------------------------------------------------------------------------
# 4. Insert a new line item row into the Monday block at row 25 (pushes the
#    Monday TOTAL row, and every row/mergeCell below it, down by one row)
# ---------------------------------------------------------------------------
$ws.Rows.Item(25).Insert()

# Copy the zebra-stripe formatting from row 23 (same style group) onto the
# freshly inserted row 25.
$ws.Range("A23:H23").Copy()
$ws.Range("A25:H25").PasteSpecial(-4122)

$ws.Range("A25").Value = "Point 11"
$ws.Range("B25").Value = "PLA-HDIG"
$ws.Range("C25").Value = "Inst"
$ws.Range("D25").Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Range("E25").Value = "EA"
$ws.Range("F25").Value = 1
$ws.Range("H25").Value = 0

# Monday TOTAL row (was row 25, now shifted to row 26)
$ws.Range("H26").Value = 0

# ---------------------------------------------------------------------------
# 5. Tuesday (07/15/2025) block — now shifted down by one row (was 30-36,
#    now 31-37)
# ---------------------------------------------------------------------------
$ws.Range("H31").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("H37").Value = 0

# ---------------------------------------------------------------------------
# 6. Wednesday (07/16/2025) block — now shifted down by one row (was 41-43,
#    now 42-44)
# ---------------------------------------------------------------------------
$ws.Range("H42").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("H44").Value = 0

# ---------------------------------------------------------------------------
# 7. Thursday (07/17/2025) block — now shifted down by one row (was 48-52,
#    now 49-53)
# ---------------------------------------------------------------------------
$ws.Range("H49").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("H53").Value = 0

Write-Output "Edit complete"
